$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column to the attendance table (Table1): "20-may"
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()
$lo.HeaderRowRange.Item(1, 5).Value = "20-may"

# Copy the header style (date number format) from the previous date column
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

# Mark attendance ("x") for the 20-may session for the students that attended
$ws.Range("E6").Value = "x"
$ws.Range("E9").Value = "x"
$ws.Range("E16").Value = "x"
$ws.Range("E20").Value = "x"
$ws.Range("E22").Value = "x"

# Update the active selection / scroll position
$ws.Range("E7").Select() | Out-Null
